$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(1)

# Remove the outline on the rounded-rectangle background shape.
$shp.Line.Visible = $false

# Update the heading text: pluralize "Religion" -> "Religions", tidy up the
# trailing comma on "Theologies", and add a new line for "and Ethics".
$nl = [char]13
$tr = $shp.TextFrame.TextRange
$tr.Text = "Worldviews, " + $nl + "Religions, " + $nl + "Theologies, " + $nl + "and Ethics"
